# Each two-digit-by-two-digit multiplication prompt in the worksheet's
# tables is replaced with a freshly generated one. Every "old" prompt is
# unique in the document, so a simple whole-document Find & Replace
# (MatchWholeWord, not MatchWildcards) targets exactly the intended cell.
$d = $word.ActiveDocument

$d.Content.Find.Execute("65×23=", $true, $false, $false, $false, $false, $true, 1, $false, "84×67=", 2) | Out-Null
$d.Content.Find.Execute("88×26=", $true, $false, $false, $false, $false, $true, 1, $false, "39×97=", 2) | Out-Null
$d.Content.Find.Execute("69×29=", $true, $false, $false, $false, $false, $true, 1, $false, "17×95=", 2) | Out-Null
$d.Content.Find.Execute("63×95=", $true, $false, $false, $false, $false, $true, 1, $false, "94×29=", 2) | Out-Null
$d.Content.Find.Execute("81×53=", $true, $false, $false, $false, $false, $true, 1, $false, "42×71=", 2) | Out-Null
$d.Content.Find.Execute("19×53=", $true, $false, $false, $false, $false, $true, 1, $false, "88×21=", 2) | Out-Null
$d.Content.Find.Execute("80×75=", $true, $false, $false, $false, $false, $true, 1, $false, "47×16=", 2) | Out-Null
$d.Content.Find.Execute("36×61=", $true, $false, $false, $false, $false, $true, 1, $false, "31×77=", 2) | Out-Null
$d.Content.Find.Execute("63×67=", $true, $false, $false, $false, $false, $true, 1, $false, "14×87=", 2) | Out-Null
$d.Content.Find.Execute("35×95=", $true, $false, $false, $false, $false, $true, 1, $false, "65×11=", 2) | Out-Null
$d.Content.Find.Execute("51×19=", $true, $false, $false, $false, $false, $true, 1, $false, "90×91=", 2) | Out-Null
$d.Content.Find.Execute("58×16=", $true, $false, $false, $false, $false, $true, 1, $false, "68×58=", 2) | Out-Null
$d.Content.Find.Execute("55×53=", $true, $false, $false, $false, $false, $true, 1, $false, "83×92=", 2) | Out-Null
$d.Content.Find.Execute("11×89=", $true, $false, $false, $false, $false, $true, 1, $false, "37×65=", 2) | Out-Null
$d.Content.Find.Execute("39×77=", $true, $false, $false, $false, $false, $true, 1, $false, "57×39=", 2) | Out-Null
$d.Content.Find.Execute("73×37=", $true, $false, $false, $false, $false, $true, 1, $false, "11×53=", 2) | Out-Null
$d.Content.Find.Execute("16×56=", $true, $false, $false, $false, $false, $true, 1, $false, "54×80=", 2) | Out-Null
$d.Content.Find.Execute("44×97=", $true, $false, $false, $false, $false, $true, 1, $false, "22×14=", 2) | Out-Null
$d.Content.Find.Execute("66×72=", $true, $false, $false, $false, $false, $true, 1, $false, "11×36=", 2) | Out-Null
$d.Content.Find.Execute("91×13=", $true, $false, $false, $false, $false, $true, 1, $false, "56×38=", 2) | Out-Null
$d.Content.Find.Execute("13×39=", $true, $false, $false, $false, $false, $true, 1, $false, "62×52=", 2) | Out-Null
$d.Content.Find.Execute("84×75=", $true, $false, $false, $false, $false, $true, 1, $false, "29×86=", 2) | Out-Null
$d.Content.Find.Execute("40×18=", $true, $false, $false, $false, $false, $true, 1, $false, "90×55=", 2) | Out-Null
$d.Content.Find.Execute("65×97=", $true, $false, $false, $false, $false, $true, 1, $false, "46×19=", 2) | Out-Null
$d.Content.Find.Execute("26×23=", $true, $false, $false, $false, $false, $true, 1, $false, "92×50=", 2) | Out-Null
